# Updated cryptos list on Sat Apr 20 13:32:02 UTC 2024 with GitHub Actions
# Applies updated price/volume(1h) figures, plus a Polygon/InternetComputer row swap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep a Text format so numeric-looking strings (e.g. "0.517",
# "63.924.16") are not auto-converted to numbers by Excel.
$changedCells = @(
    "D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "E7", "D8", "E8", "D9", "E9", "E10", "D11",
    "E11", "D12", "E12", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "E18", "D19",
    "E19", "D20", "E20", "D21", "E21", "B22", "C22", "D22", "E22", "B23", "C23", "D23", "E23", "E24",
    "D25", "E25", "E26", "E27", "D28", "E28", "D29", "E29", "E30", "D31", "E31", "E32", "D33", "E33",
    "D34", "E34", "D35", "E35", "D36", "E36", "D37", "E37", "D38", "E38", "D39", "E39", "D40", "E40",
    "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "E48", "D49",
    "E49", "D50", "E50", "E51"
)
foreach ($cellAddr in $changedCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.924.16"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "3.063.47"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "558.65"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").Value = "142.66"
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.062.67"
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("D9").Value = "0.517"
$ws.Range("E9").Value = "  +3.69%  "
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("D11").Value = "6.19"
$ws.Range("E11").Value = "  -2.54%  "
$ws.Range("D12").Value = "0.480"
$ws.Range("E12").Value = "  +2.25%  "
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("D14").Value = "35.25"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "3.563.28"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("D16").Value = "63.893.50"
$ws.Range("E16").Value = "  -1.17%  "
$ws.Range("D17").Value = "3.055.86"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").Value = "6.78"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "489.54"
$ws.Range("E20").Value = "  +2.42%  "
$ws.Range("D21").Value = "14.27"
$ws.Range("E21").Value = "  +3.88%  "
$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").Value = "14.72"
$ws.Range("E22").Value = "  +9.61%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").Value = "0.686"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("D25").Value = "82.81"
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("D28").Value = "8.14"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").Value = "2.05"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").Value = "26.52"
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").Value = "2.53"
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("D34").Value = "5.70"
$ws.Range("E34").Value = "  +1.71%  "
$ws.Range("D35").Value = "6.23"
$ws.Range("E35").Value = "  +1.40%  "
$ws.Range("D36").Value = "55.36"
$ws.Range("E36").Value = "  +0.82%  "
$ws.Range("D37").Value = "0.0410"
$ws.Range("E37").Value = "  +0.59%  "
$ws.Range("D38").Value = "444.23"
$ws.Range("E38").Value = "  -4.11%  "
$ws.Range("D39").Value = "0.0816"
$ws.Range("E39").Value = "  -2.00%  "
$ws.Range("D40").Value = "3.031.63"
$ws.Range("E40").Value = "  +1.89%  "
$ws.Range("D41").Value = "2.80"
$ws.Range("E41").Value = "  -5.65%  "
$ws.Range("D42").Value = "8.34"
$ws.Range("E42").Value = "  +1.19%  "
$ws.Range("D43").Value = "0.117"
$ws.Range("E43").Value = "  +1.82%  "
$ws.Range("D44").Value = "0.274"
$ws.Range("E44").Value = "  +5.92%  "
$ws.Range("D45").Value = "27.75"
$ws.Range("E45").Value = "  -1.16%  "
$ws.Range("D46").Value = "2.25"
$ws.Range("E46").Value = "  +4.71%  "
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("D49").Value = "118.26"
$ws.Range("E49").Value = "  +1.40%  "
$ws.Range("D50").Value = "0.0₃0517"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("E51").Value = "  +2.85%  "
